$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.385.36"
$ws.Range("E2").Value = "  +1.59%  "

$ws.Range("D3").Value = "2.158.04"
$ws.Range("E3").Value = "  +3.13%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'227.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.53%  "

$ws.Range("D6").Value = "'0.622"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.85%  "

$ws.Range("D7").Value = "'64.12"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.70%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  +2.86%  "

$ws.Range("D10").Value = "'0.0860"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.66%  "

$ws.Range("E11").Value = "  -0.15%  "

$ws.Range("D12").Value = "'16.00"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.30%  "

$ws.Range("D13").Value = "2.477.52"
$ws.Range("E13").Value = "  +3.09%  "

$ws.Range("D14").Value = "'22.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.85%  "

$ws.Range("D15").Value = "'0.813"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.92%  "

$ws.Range("E16").Value = "  +0.92%  "

$ws.Range("D17").Value = "2.160.33"
$ws.Range("E17").Value = "  +3.42%  "

$ws.Range("D18").Value = "39.359.83"
$ws.Range("E18").Value = "  +1.71%  "

$ws.Range("E19").Value = "  -0.03%  "

$ws.Range("E20").Value = "  +0.48%  "

$ws.Range("E21").Value = "  +1.38%  "

$ws.Range("D22").Value = "'231.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.54%  "

$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("E24").Value = "  +5.87%  "

$ws.Range("D25").Value = "'2.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.66%  "

$ws.Range("D26").Value = "'172.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.34%  "

$ws.Range("D27").Value = "'9.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("E28").Value = "  +0.77%  "

$ws.Range("E29").Value = "  +2.97%  "

$ws.Range("E30").Value = "  -1.24%  "

$ws.Range("D31").Value = "'2.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.66%  "

$ws.Range("E32").Value = "  +0.89%  "

$ws.Range("E33").Value = "  +2.04%  "

$ws.Range("E34").Value = "  +0.23%  "

$ws.Range("E35").Value = "  +8.99%  "

$ws.Range("E36").Value = "  -0.36%  "

$ws.Range("E37").Value = "  +0.32%  "

$ws.Range("E38").Value = "  -0.15%  "

$ws.Range("E39").Value = "  +0.22%  "

$ws.Range("D40").Value = "'103.87"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.87%  "

$ws.Range("D41").Value = "'0.0230"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.71%  "

$ws.Range("D42").Value = "'17.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.19%  "

$ws.Range("D43").Value = "1.540.99"
$ws.Range("E43").Value = "  +0.33%  "

$ws.Range("E44").Value = "  +3.80%  "

$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'7.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.39%  "

$ws.Range("B46").Value = "HuobiToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D46").Value = "'2.82"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.61%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.0925"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.67%  "

$ws.Range("E48").Value = "  +5.67%  "

$ws.Range("E49").Value = "  +2.78%  "

$ws.Range("D50").Value = "2.361.32"
$ws.Range("E50").Value = "  +3.11%  "

$ws.Range("E51").Value = "  +0.17%  "
